# Update the CDA Logical model metadata (StructureDefinition-CS) for ST.r2b.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Bump Version and Date values (row 3 = Version, row 8 = Date).
$ws.Cells.Item(3, 2).Value2 = "2.0.1-sd-202510-matchbox-patch"
$ws.Cells.Item(8, 2).Value2 = "2025-10-29T22:15:57+01:00"

# 2. Insert a new "Jurisdiction" property row right after "Contact" (row 10),
#    pushing "Description" and everything below it down by one row.
$ws.Rows.Item(11).Insert()

# Match the formatting of the surrounding property rows (border/wrap style)
# instead of the bare default style a blank Insert() would leave behind.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(11, 1).Value2 = "Jurisdiction"
$ws.Cells.Item(11, 2).Value2 = ""
